# Update the date line and the 25 division problems in the single table.
$d = $word.ActiveDocument

# --- Title paragraph date update ---
$d.Paragraphs.Item(1).Range.Text = "2023-08-31 Thursday"

# --- Table cell updates ---
# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-indexed)
# contain problems. Replacements are applied by explicit (row, col)
# addressing so that overlapping old/new values never collide.
$t = $d.Tables.Item(1)

$rows = @(1, 5, 9, 13, 17)

$values = @(
    @("79÷2=", "79÷8=", "19÷6=", "40÷8=", "35÷9="),
    @("85÷2=", "71÷5=", "50÷9=", "13÷4=", "98÷2="),
    @("60÷5=", "26÷5=", "38÷4=", "43÷4=", "15÷3="),
    @("39÷3=", "63÷9=", "21÷8=", "33÷9=", "33÷5="),
    @("34÷9=", "81÷9=", "76÷5=", "80÷7=", "94÷8=")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $rowValues = $values[$i]
    for ($c = 1; $c -le 5; $c++) {
        $t.Cell($r, $c).Range.Text = $rowValues[$c - 1]
    }
}
